# Auto-generated edit script applying scheduled-runner data refresh to Lamia_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4589.25
$ws.Range("J17").Value = 4589.25
$ws.Range("L17").Value = 13767.75
$ws.Range("N17").Value = -14103.75
$ws.Range("H18").Value = 928.8
$ws.Range("I18").Value = 928.8
$ws.Range("K18").Value = 928.8
$ws.Range("M18").Value = -644.8
$ws.Range("H33").Value = 364.42307
$ws.Range("I33").Value = 360
$ws.Range("K33").Value = 360
$ws.Range("M33").Value = -131
$ws.Range("H43").Value = 6389.5
$ws.Range("I43").Value = 2499.5
$ws.Range("J43").Value = 7362
$ws.Range("K43").Value = 2499.5
$ws.Range("L43").Value = 7362
$ws.Range("M43").Value = -2430.5
$ws.Range("N43").Value = -7500
$ws.Range("H62").Value = 6060.722
$ws.Range("I62").Value = 3014.6667
$ws.Range("J62").Value = 7583.75
$ws.Range("K62").Value = 3014.6667
$ws.Range("L62").Value = 7583.75
$ws.Range("M62").Value = -2390.6667
$ws.Range("N62").Value = -8831.75
$ws.Range("H65").Value = 6060.722
$ws.Range("I65").Value = 3014.6667
$ws.Range("J65").Value = 7583.75
$ws.Range("K65").Value = 15073.3335
$ws.Range("L65").Value = 37918.75
$ws.Range("M65").Value = -11953.3335
$ws.Range("N65").Value = -44158.75
$ws.Range("H76").Value = 7691.6665
$ws.Range("I76").Value = 4749.5
$ws.Range("K76").Value = 4749.5
$ws.Range("M76").Value = -4434.5
$ws.Range("H79").Value = 7691.6665
$ws.Range("I79").Value = 4749.5
$ws.Range("K79").Value = 4749.5
$ws.Range("M79").Value = -3657.5
$ws.Range("H88").Value = 205413.53
$ws.Range("I88").Value = 603580
$ws.Range("K88").Value = 603580
$ws.Range("M88").Value = -603174
$ws.Range("H91").Value = 205413.53
$ws.Range("I91").Value = 603580
$ws.Range("K91").Value = 603580
$ws.Range("M91").Value = -602176
$ws.Range("H98").Value = 437829.78
$ws.Range("I98").Value = 1011
$ws.Range("K98").Value = 1011
$ws.Range("M98").Value = 487
$ws.Range("H100").Value = 7664.727
$ws.Range("I100").Value = 7356.3335
$ws.Range("J100").Value = 8034.8
$ws.Range("K100").Value = 7356.3335
$ws.Range("L100").Value = 8034.8
$ws.Range("M100").Value = -6815.3335
$ws.Range("N100").Value = -9116.799999999999
$ws.Range("H106").Value = 11180.077
$ws.Range("I106").Value = 7088.1055
$ws.Range("K106").Value = 7088.1055
$ws.Range("M106").Value = -6457.1055
$ws.Range("H122").Value = 437829.78
$ws.Range("I122").Value = 1011
$ws.Range("K122").Value = 3033
$ws.Range("M122").Value = -583
$ws.Range("H125").Value = 1882.6129
$ws.Range("I125").Value = 683.9
$ws.Range("J125").Value = 2453.4285
$ws.Range("K125").Value = 6155.099999999999
$ws.Range("L125").Value = 22080.8565
$ws.Range("M125").Value = -3695.099999999999
$ws.Range("N125").Value = -27000.8565
$ws.Range("H132").Value = 1556.6
$ws.Range("I132").Value = 1412.7407
$ws.Range("J132").Value = 2851.3333
$ws.Range("K132").Value = 4238.2221
$ws.Range("L132").Value = 8553.999899999999
$ws.Range("M132").Value = -1708.2221
$ws.Range("N132").Value = -13613.9999
$ws.Range("H137").Value = 20410704
$ws.Range("I137").Value = 52633656
$ws.Range("J137").Value = 2836.2334
$ws.Range("K137").Value = 157900968
$ws.Range("L137").Value = 8508.700199999999
$ws.Range("M137").Value = -157898418
$ws.Range("N137").Value = -13608.7002

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16
$ws.Range("H32").Value = 3171.3052
$ws.Range("I32").Value = 3171.3052
$ws.Range("K32").Value = 3171.3052
$ws.Range("M32").Value = -2884.3052
$ws.Range("H45").Value = 1928.9166
$ws.Range("I45").Value = 1572.1111
$ws.Range("J45").Value = 2999.3333
$ws.Range("K45").Value = 1572.1111
$ws.Range("L45").Value = 2999.3333
$ws.Range("M45").Value = -1195.1111
$ws.Range("N45").Value = -3753.3333
$ws.Range("H74").Value = 5954503.5
$ws.Range("I74").Value = 6537391
$ws.Range("J74").Value = 9051.4
$ws.Range("K74").Value = 6537391
$ws.Range("L74").Value = 9051.4
$ws.Range("M74").Value = -6536517
$ws.Range("N74").Value = -10799.4
$ws.Range("H77").Value = 5954503.5
$ws.Range("I77").Value = 6537391
$ws.Range("J77").Value = 9051.4
$ws.Range("K77").Value = 32686955
$ws.Range("L77").Value = 45257
$ws.Range("M77").Value = -32682587
$ws.Range("N77").Value = -53993
$ws.Range("H122").Value = 76927670
$ws.Range("I122").Value = 2928
$ws.Range("K122").Value = 8784
$ws.Range("M122").Value = -6334
$ws.Range("H132").Value = 3489.673
$ws.Range("I132").Value = 3107.9092
$ws.Range("K132").Value = 9323.7276
$ws.Range("M132").Value = -6793.7276
$ws.Range("H134").Value = 75429
$ws.Range("J134").Value = 75429
$ws.Range("L134").Value = 75429
$ws.Range("N134").Value = -85569

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 19500
$ws.Range("J32").Value = 19500
$ws.Range("L32").Value = 19500
$ws.Range("N32").Value = -20268
$ws.Range("H86").Value = 2464.75
$ws.Range("I86").Value = 2080.9644
$ws.Range("K86").Value = 2080.9644
$ws.Range("M86").Value = -957.9643999999998
$ws.Range("H89").Value = 2464.75
$ws.Range("I89").Value = 2080.9644
$ws.Range("K89").Value = 10404.822
$ws.Range("M89").Value = -4788.822

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35640.883
$ws.Range("I31").Value = 3468
$ws.Range("J31").Value = 102911.45
$ws.Range("K31").Value = 3468
$ws.Range("L31").Value = 102911.45
$ws.Range("M31").Value = -3173
$ws.Range("N31").Value = -103501.45
$ws.Range("H34").Value = 35640.883
$ws.Range("I34").Value = 3468
$ws.Range("J34").Value = 102911.45
$ws.Range("K34").Value = 3468
$ws.Range("L34").Value = 102911.45
$ws.Range("M34").Value = -3266
$ws.Range("N34").Value = -103315.45
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = $null
$ws.Range("H58").Value = 4503.696
$ws.Range("I58").Value = 2091.375
$ws.Range("J58").Value = 10017.571
$ws.Range("K58").Value = 2091.375
$ws.Range("L58").Value = 10017.571
$ws.Range("M58").Value = -1888.375
$ws.Range("N58").Value = -10423.571
$ws.Range("H134").Value = 8096.722
$ws.Range("I134").Value = 7750.8
$ws.Range("J134").Value = 8882.909
$ws.Range("K134").Value = 23252.4
$ws.Range("L134").Value = 26648.727
$ws.Range("M134").Value = -20717.4
$ws.Range("N134").Value = -31718.727
$ws.Range("H136").Value = 4503.696
$ws.Range("I136").Value = 2091.375
$ws.Range("J136").Value = 10017.571
$ws.Range("K136").Value = 6274.125
$ws.Range("L136").Value = 30052.713
$ws.Range("M136").Value = -3724.125
$ws.Range("N136").Value = -35152.713

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2928.25
$ws.Range("I140").Value = 2354.5
$ws.Range("K140").Value = 7063.5
$ws.Range("M140").Value = -1883.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5949.523
$ws.Range("I40").Value = 5876.825
$ws.Range("K40").Value = 5876.825
$ws.Range("M40").Value = -5740.825
$ws.Range("H46").Value = 4854.727
$ws.Range("I46").Value = 3600
$ws.Range("J46").Value = 5571.7144
$ws.Range("K46").Value = 3600
$ws.Range("L46").Value = 5571.7144
$ws.Range("M46").Value = -3412
$ws.Range("N46").Value = -5947.7144
$ws.Range("H122").Value = 129608.75
$ws.Range("I122").Value = 164143.16
$ws.Range("J122").Value = 6271.5713
$ws.Range("K122").Value = 492429.48
$ws.Range("L122").Value = 18814.7139
$ws.Range("M122").Value = -489979.48
$ws.Range("N122").Value = -23714.7139
$ws.Range("H132").Value = 7856.579
$ws.Range("I132").Value = 6734.7144
$ws.Range("J132").Value = 10997.8
$ws.Range("K132").Value = 20204.1432
$ws.Range("L132").Value = 32993.39999999999
$ws.Range("M132").Value = -17674.1432
$ws.Range("N132").Value = -38053.39999999999
$ws.Range("H136").Value = 5293.645
$ws.Range("I136").Value = 4194.3335
$ws.Range("J136").Value = 10390.454
$ws.Range("K136").Value = 12583.0005
$ws.Range("L136").Value = 31171.362
$ws.Range("M136").Value = -10033.0005
$ws.Range("N136").Value = -36271.362

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 459.8
$ws.Range("I2").Value = 459.8
$ws.Range("K2").Value = 459.8
$ws.Range("M2").Value = -347.8
$ws.Range("H62").Value = 14040.533
$ws.Range("I62").Value = 7899.3335
$ws.Range("J62").Value = 15575.833
$ws.Range("K62").Value = 7899.3335
$ws.Range("L62").Value = 15575.833
$ws.Range("M62").Value = -7275.3335
$ws.Range("N62").Value = -16823.833
$ws.Range("H65").Value = 14040.533
$ws.Range("I65").Value = 7899.3335
$ws.Range("J65").Value = 15575.833
$ws.Range("K65").Value = 39496.6675
$ws.Range("L65").Value = 77879.16500000001
$ws.Range("M65").Value = -36376.6675
$ws.Range("N65").Value = -84119.16500000001
$ws.Range("H132").Value = 3276.606
$ws.Range("I132").Value = 1713.762
$ws.Range("J132").Value = 6011.5835
$ws.Range("K132").Value = 5141.286
$ws.Range("L132").Value = 18034.7505
$ws.Range("M132").Value = -2611.286
$ws.Range("N132").Value = -23094.7505
$ws.Range("H136").Value = 2293.5952
$ws.Range("I136").Value = 1495.3513
$ws.Range("J136").Value = 8200.6
$ws.Range("K136").Value = 4486.0539
$ws.Range("L136").Value = 24601.8
$ws.Range("M136").Value = -1936.0539
$ws.Range("N136").Value = -29701.8
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

Write-Output "Applied 253 cell updates across 7 sheets"
